$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rHigh = $ws.Range("F111:F137")
$rHigh.VerticalAlignment = -4160
$rHigh.NumberFormat = "#,##0"

$rLow = $ws.Range("F2:F110")
$rLow.VerticalAlignment = -4160
$rLow.NumberFormat = "0"

$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 18
$ws.Range("F5").Value = 34
$ws.Range("F6").Value = 17
$ws.Range("F7").Value = 16
$ws.Range("F8").Value = 6
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 15
$ws.Range("F11").Value = 12
$ws.Range("F12").Value = 17
$ws.Range("F13").Value = 21
$ws.Range("F14").Value = 17
$ws.Range("F15").Value = 18
$ws.Range("F16").Value = 7
$ws.Range("F17").Value = 6
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 5
$ws.Range("F20").Value = 5
$ws.Range("F21").Value = 8
$ws.Range("F22").Value = 6
$ws.Range("F23").Value = 1
$ws.Range("F24").Value = 5
$ws.Range("F25").Value = 9
$ws.Range("F26").Value = 3
$ws.Range("F27").Value = 5
$ws.Range("F28").Value = 5
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 20
$ws.Range("F31").Value = 4
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 13
$ws.Range("F34").Value = 22
$ws.Range("F35").Value = 8
$ws.Range("F36").Value = 9
$ws.Range("F37").Value = 7
$ws.Range("F38").Value = 27
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 10
$ws.Range("F41").Value = 10
$ws.Range("F42").Value = 17
$ws.Range("F43").Value = 6
$ws.Range("F44").Value = 30
$ws.Range("F45").Value = 21
$ws.Range("F46").Value = 30
$ws.Range("F47").Value = 4
$ws.Range("F48").Value = 6
$ws.Range("F49").Value = 10
$ws.Range("F50").Value = 0
$ws.Range("F51").Value = 14
$ws.Range("F52").Value = 2
$ws.Range("F53").Value = 42
$ws.Range("F54").Value = 26
$ws.Range("F55").Value = 16
$ws.Range("F56").Value = 6
$ws.Range("F57").Value = 7
$ws.Range("F58").Value = 0
$ws.Range("F59").Value = 0
$ws.Range("F60").Value = 4
$ws.Range("F61").Value = -1
$ws.Range("F62").Value = 27
$ws.Range("F63").Value = 36
$ws.Range("F64").Value = 27
$ws.Range("F65").Value = 40
$ws.Range("F66").Value = 14
$ws.Range("F67").Value = 10
$ws.Range("F68").Value = 0
$ws.Range("F69").Value = 14
$ws.Range("F70").Value = 10
$ws.Range("F71").Value = 0
$ws.Range("F72").Value = 0
$ws.Range("F73").Value = 10
$ws.Range("F74").Value = 22
$ws.Range("F75").Value = 16
$ws.Range("F76").Value = 17
$ws.Range("F77").Value = 15
$ws.Range("F78").Value = 2
$ws.Range("F79").Value = 9
$ws.Range("F80").Value = 1
$ws.Range("F81").Value = 27
$ws.Range("F82").Value = 12
$ws.Range("F83").Value = 21
$ws.Range("F84").Value = 18
$ws.Range("F85").Value = 21
$ws.Range("F86").Value = 4
$ws.Range("F87").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("F89").Value = 8
$ws.Range("F90").Value = 9
$ws.Range("F91").Value = 36
$ws.Range("F92").Value = 55
$ws.Range("F93").Value = 16
$ws.Range("F94").Value = 9
$ws.Range("F95").Value = 5
$ws.Range("F96").Value = 13
$ws.Range("F97").Value = 0
$ws.Range("F98").Value = 7
$ws.Range("F99").Value = 9
$ws.Range("F100").Value = 37
$ws.Range("F101").Value = 21
$ws.Range("F102").Value = 15
$ws.Range("F103").Value = 4
$ws.Range("F104").Value = 4
$ws.Range("F105").Value = 27
$ws.Range("F106").Value = 34
$ws.Range("F107").Value = 14
$ws.Range("F108").Value = 13
$ws.Range("F109").Value = 3
$ws.Range("F110").Value = 1
$ws.Range("F111").Value = 0
$ws.Range("F112").Value = 1
$ws.Range("F113").Value = 32
$ws.Range("F114").Value = 0
$ws.Range("F115").Value = 0
$ws.Range("F116").Value = 8
$ws.Range("F117").Value = 0
$ws.Range("F118").Value = 27
$ws.Range("F119").Value = 22
$ws.Range("F120").Value = 30
$ws.Range("F121").Value = 14
$ws.Range("F122").Value = 19
$ws.Range("F123").Value = 10
$ws.Range("F124").Value = 9
$ws.Range("F125").Value = 23
$ws.Range("F126").Value = 25
$ws.Range("F127").Value = 11
$ws.Range("F128").Value = 0
$ws.Range("F129").Value = 8
$ws.Range("F130").Value = 8
$ws.Range("F131").Value = 8
$ws.Range("F132").Value = 3
$ws.Range("F133").Value = 6
$ws.Range("F134").Value = 0
$ws.Range("F135").Value = 0
$ws.Range("F136").Value = 7
$ws.Range("F137").Value = 3
